$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host $ws.Name
Write-Host $ws.Range("B3").Value
Write-Host $ws.Range("B4").Value
